# Updated cryptos list values (price + volume) to match the latest scrape.
# Also swaps the WhiteBITCoin / Stacks row positions (ranks 46/47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.212.29"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "'3.648.49"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'240.52"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  +12.49%  "
$ws.Range("D7").Value = "'660.32"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").Value = "'0.421"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'3.643.92"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "'44.53"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "'0.203"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'6.60"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").Value = "'4.326.80"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "'0.0000267"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").Value = "'96.059.92"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "'8.79"
$ws.Range("E18").Value = "  +13.46%  "
$ws.Range("D19").Value = "'3.641.60"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'12.63"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'18.16"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "'520.11"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "'3.42"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'6.83"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("D27").Value = "'102.04"
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("D28").Value = "'12.88"
$ws.Range("E28").Value = "  -1.95%  "
$ws.Range("E29").Value = "  +7.26%  "
$ws.Range("D30").Value = "'3.01"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("D31").Value = "'12.04"
$ws.Range("E31").Value = "  +4.25%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  +9.85%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'32.49"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").Value = "'0.580"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "'617.09"
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "'8.67"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").Value = "'42.87"
$ws.Range("E40").Value = "  +29.64%  "
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").Value = "'0.943"
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").Value = "'1.92"
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'6.11"
$ws.Range("E45").Value = "  +5.25%  "
$ws.Range("D46").Value = "'0.0449"
$ws.Range("E46").Value = "  +3.98%  "
$ws.Range("D47").Value = "'0.425"
$ws.Range("E47").Value = "  +16.84%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "'23.58"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.27"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'8.44"
$ws.Range("E50").Value = "  +1.22%  "
$ws.Range("E51").Value = "  -0.09%  "
